# Update the "Avverkningsanmälningar" worksheet:
#  - Column C ("Förändrad") is bumped from 2023-09-21 (45190) to 2023-09-23 (45192)
#    for every existing data row (rows 2-349).
#  - Row 349 picks up an explicit row height (15pt / customHeight).
#  - A new data row (350) is appended for case "A 44906-2023".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" (changed) date column for all existing data rows.
for ($i = 2; $i -le 349; $i++) {
    $ws.Cells.Item($i, 3).Value = 45192
}

# Row 349 gets an explicit custom row height in the new workbook.
$ws.Rows.Item(349).RowHeight = 15

# Append the new row (350) with the new case data.
$ws.Cells.Item(350, 1).Value = "A 44906-2023"

$ws.Cells.Item(350, 2).Value = 45190
$ws.Cells.Item(350, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(350, 3).Value = 45192
$ws.Cells.Item(350, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(350, 4).Value = "VÄSTMANLANDS LÄN"
$ws.Cells.Item(350, 5).Value = "SKINNSKATTEBERG"
$ws.Cells.Item(350, 6).Value = "Sveaskog"

$ws.Cells.Item(350, 7).Value = 1.8
$ws.Cells.Item(350, 8).Value = 0
$ws.Cells.Item(350, 9).Value = 0
$ws.Cells.Item(350, 10).Value = 0
$ws.Cells.Item(350, 11).Value = 0
$ws.Cells.Item(350, 12).Value = 0
$ws.Cells.Item(350, 13).Value = 0
$ws.Cells.Item(350, 14).Value = 0
$ws.Cells.Item(350, 15).Value = 0
$ws.Cells.Item(350, 16).Value = 0
$ws.Cells.Item(350, 17).Value = 0

# R350 mirrors the (empty, wrap-text styled) "Artnamn" cell used throughout
# column R for rows without species names.
$ws.Cells.Item(350, 18).WrapText = $true
